$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1K_ADAMTS19_01", 0, 47),
    @("1K_ADAMTS19_02", 0, 190),
    @("1K_ADAMTS19_03", 0, 133),
    @("1K_ADAMTS19_04", 0, 98),
    @("1K_ADAMTS19_04", 0, 15),
    @("1K_ADAMTS19_10", 0, 24),
    @("1K_ADAMTS19_10", 1, 16),
    @("1K_ADAMTS19_10", 2, 14),
    @("1K_ADAMTS19_10", 0, 27),
    @("1K_ADAMTS19_10", 1, 35),
    @("1K_ADAMTS19_10", 2, 11),
    @("1K_ADAMTS19_10", 3, 19),
    @("1K_ADAMTS19_10", 4, 24),
    @("1K_ADAMTS19_10", 5, 25),
    @("1K_SCR_01", 0, 16),
    @("1K_SCR_01", 1, 57),
    @("1K_SCR_01", 2, 18),
    @("1K_SCR_01", 3, 19),
    @("1K_SCR_01", 4, 31),
    @("1K_SCR_01", 5, 33),
    @("1K_SCR_01", 6, 56),
    @("1K_SCR_01", 7, 23),
    @("1K_SCR_01", 8, 20),
    @("1K_SCR_01", 9, 13),
    @("1K_SCR_01", 10, 50),
    @("1K_SCR_01", 11, 12),
    @("1K_SCR_01", 12, 21),
    @("1K_SCR_01", 13, 30),
    @("1K_SCR_01", 14, 24),
    @("1K_SCR_01", 15, 42),
    @("1K_SCR_01", 16, 14),
    @("1K_SCR_01", 0, 22),
    @("1K_SCR_03", 0, 61),
    @("1K_SCR_03", 1, 340),
    @("1K_SCR_05", 0, 30),
    @("1K_SCR_05", 0, 15),
    @("1K_TAX1BP1_01", 0, 103),
    @("1K_TAX1BP1_02", 0, 15),
    @("1K_TAX1BP1_02", 1, 129),
    @("1K_TAX1BP1_02", 2, 44),
    @("1K_TAX1BP1_02", 3, 127),
    @("1K_TAX1BP1_02", 4, 12),
    @("1K_TAX1BP1_02", 5, 15),
    @("1K_TAX1BP1_02", 6, 16),
    @("1K_TAX1BP1_02", 7, 22),
    @("1K_TAX1BP1_02", 8, 33),
    @("1K_TAX1BP1_02", 9, 23),
    @("1K_TAX1BP1_04", 0, 11),
    @("1K_TAX1BP1_04", 1, 24),
    @("1K_TAX1BP1_04", 2, 30),
    @("1K_TAX1BP1_04", 3, 12),
    @("1K_TAX1BP1_04", 4, 21),
    @("1K_TAX1BP1_04", 5, 22),
    @("1K_TAX1BP1_04", 6, 37),
    @("1K_TAX1BP1_04", 7, 13),
    @("1K_TAX1BP1_04", 8, 28),
    @("1K_TAX1BP1_04", 9, 154),
    @("1K_TAX1BP1_04", 10, 47),
    @("1K_TAX1BP1_04", 11, 54),
    @("1K_TAX1BP1_04", 12, 23),
    @("1K_TAX1BP1_04", 13, 174),
    @("1K_TAX1BP1_04", 14, 91),
    @("1K_TAX1BP1_04", 15, 28),
    @("1K_TAX1BP1_04", 16, 21),
    @("1K_TAX1BP1_05", 0, 22),
    @("1K_TAX1BP1_05", 1, 16),
    @("1K_TAX1BP1_05", 2, 28),
    @("1K_TAX1BP1_05", 3, 16),
    @("1K_TAX1BP1_07", 0, 416),
    @("1K_TAX1BP1_08", 0, 63),
    @("1K_TAX1BP1_09", 0, 113),
    @("1K_TAX1BP1_10", 0, 458),
    @("WT_ADAMTS19_01", 0, 265),
    @("WT_ADAMTS19_02", 0, 87),
    @("WT_ADAMTS19_07", 0, 182),
    @("WT_ADAMTS19_10", 0, 773),
    @("WT_SCR_01", 0, 243),
    @("WT_SCR_03", 0, 16),
    @("WT_SCR_03", 1, 35),
    @("WT_SCR_03", 2, 15),
    @("WT_SCR_05", 0, 35),
    @("WT_SCR_05", 1, 14),
    @("WT_SCR_07", 0, 13),
    @("WT_SCR_08", 0, 16),
    @("WT_SCR_08", 0, 173),
    @("WT_SCR_08", 1, 16),
    @("WT_SCR_08", 2, 55),
    @("WT_SCR_10", 0, 298),
    @("WT_TAX1BP1_01", 0, 76),
    @("WT_TAX1BP1_01", 1, 663),
    @("WT_TAX1BP1_02", 0, 25),
    @("WT_TAX1BP1_03", 0, 465),
    @("WT_TAX1BP1_03", 0, 195),
    @("WT_TAX1BP1_04", 0, 1800),
    @("WT_TAX1BP1_04", 0, 606),
    @("WT_TAX1BP1_05", 0, 15),
    @("WT_TAX1BP1_06", 0, 19),
    @("WT_TAX1BP1_06", 1, 23),
    @("WT_TAX1BP1_06", 2, 27),
    @("WT_TAX1BP1_07", 0, 221),
    @("WT_TAX1BP1_08", 0, 201),
    @("WT_TAX1BP1_09", 0, 111),
    @("WT_TAX1BP1_10", 0, 13),
    @("WT_TAX1BP1_10", 1, 80),
    @("WT_TAX1BP1_10", 2, 34),
    @("WT_TAX1BP1_10", 3, 11),
    @("WT_TAX1BP1_10", 4, 22),
    @("WT_TAX1BP1_10", 5, 12),
    @("WT_TAX1BP1_10", 6, 11),
    @("WT_TAX1BP1_10", 7, 29)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
